# Fixed Bento 80 Test scripts
# Appends an " order By ... LIMIT 100" clause to each of the three
# Cypher queries stored on the "startup" sheet (columns B2, B3, B4),
# matching the author's update to the Neo4j queries used by the
# Cases / Samples / Files tabs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B2: "CasesTab" query -> gets "order By ss.study_subject_id" ---
$cell = $ws.Range("B2")
$orig = [string]$cell.Value2
if ($orig -notmatch "order By ss\.study_subject_id") {
    $cell.Value2 = $orig + "`n order By ss.study_subject_id ASC LIMIT 100 "
}

# --- B3: "SamplesTab" query -> gets "order By samp.sample_id" ---
$cell = $ws.Range("B3")
$orig = [string]$cell.Value2
if ($orig -notmatch "order By samp\.sample_id") {
    $cell.Value2 = $orig + "`n order By samp.sample_id ASC LIMIT 100"
}

# --- B4: "FilesTab" query -> gets "order By f.file_name" ---
$cell = $ws.Range("B4")
$orig = [string]$cell.Value2
if ($orig -notmatch "order By f\.file_name") {
    $cell.Value2 = $orig + "`n order By f.file_name ASC LIMIT 100"
}

# --- Selection moved from C4 to B4 (also resets the scrolled-to-A4 view) ---
$ws.Range("B4").Select()
